$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.637.33"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.597.67"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Formula = "'211.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Formula = "'0.510"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Formula = "'19.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").Formula = "'0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.821.28"
$ws.Range("D13").Value = "1.589.91"
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").Formula = "'4.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Formula = "'0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Formula = "'64.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "26.617.04"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("E18").Value = "  -2.63%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Formula = "'208.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").Formula = "'7.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Formula = "'2.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").Formula = "'8.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Formula = "'143.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Formula = "'7.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Formula = "'0.113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").Formula = "'15.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +18.90%  "
$ws.Range("D35").Value = "1.276.62"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("D40").Formula = "'0.821"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Formula = "'62.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").Value = "1.733.46"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Formula = "'89.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("D47").Formula = "'1.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Formula = "'0.0512"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("E51").Value = "  +0.28%  "
